$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the dSF (column F) values to reflect the repulled data.
# One row (61) also needed its dS0 (column E) value corrected.
$ws.Range("F7").Value = -2
$ws.Range("F11").Value = -2
$ws.Range("F13").Value = 1
$ws.Range("F17").Value = 2
$ws.Range("F23").Value = -1
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = 1
$ws.Range("F30").Value = -7
$ws.Range("F31").Value = 0
$ws.Range("F35").Value = -2
$ws.Range("F40").Value = 2
$ws.Range("F41").Value = -4
$ws.Range("F46").Value = -3
$ws.Range("F47").Value = 2
$ws.Range("F55").Value = 2
$ws.Range("F56").Value = -3
$ws.Range("F57").Value = 5
$ws.Range("F58").Value = -4
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 1
$ws.Range("F62").Value = -2
$ws.Range("F63").Value = 2
$ws.Range("F68").Value = -2
$ws.Range("F70").Value = -1
